$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the relative influence values for rows 2-5 (labels unchanged)
$ws.Range("B2").Value = 35.01803606954568
$ws.Range("B3").Value = 18.27220470233163
$ws.Range("B4").Value = 16.34857708301666
$ws.Range("B5").Value = 12.22379807659952

# Rows 6 and 7 swap labels (BottomDepth <-> SuspendedParticulateMatter) and get new values
$ws.Range("A6").Value = "SuspendedParticulateMatter"
$ws.Range("B6").Value = 4.560418877824543
$ws.Range("A7").Value = "BottomDepth"
$ws.Range("B7").Value = 3.899998206344125

# Update remaining rows
$ws.Range("B8").Value = 3.668622449248512
$ws.Range("B9").Value = 3.533720293118118
$ws.Range("B10").Value = 2.474624241971222
